$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 11.92973733333333
$ws.Cells.Item(2, 8).Value = 35.789212
$ws.Cells.Item(2, 9).Value = 0.1260232909831994
$ws.Cells.Item(2, 10).Value = 0.1260232909831994
$ws.Cells.Item(2, 13).Value = 3.795192333333334
$ws.Cells.Item(2, 14).Value = 11.385577
$ws.Cells.Item(2, 15).Value = 0.01044213755712683
$ws.Cells.Item(2, 16).Value = 0.01044213755712683
$ws.Cells.Item(2, 17).Value = 45.27564766614712
$ws.Cells.Item(2, 18).Value = 407.480828995324
$ws.Cells.Item(2, 19).Value = 0.00131595253984839
$ws.Cells.Item(2, 20).Value = 0.00131595253984839
$ws.Cells.Item(3, 7).Value = 11.92973733333333
$ws.Cells.Item(3, 8).Value = 35.789212
$ws.Cells.Item(3, 9).Value = 0.1260232909831994
$ws.Cells.Item(3, 10).Value = 0.1260232909831994
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.6696287328350964
$ws.Cells.Item(3, 16).Value = 0.6696287328350964
$ws.Cells.Item(3, 17).Value = 2903.41650922595
$ws.Cells.Item(3, 18).Value = 26130.74858303355
$ws.Cells.Item(3, 19).Value = 0.08438881664878847
$ws.Cells.Item(3, 20).Value = 0.08438881664878847
$ws.Cells.Item(4, 7).Value = 11.92973733333333
$ws.Cells.Item(4, 8).Value = 35.789212
$ws.Cells.Item(4, 9).Value = 0.1260232909831994
$ws.Cells.Item(4, 10).Value = 0.1260232909831994
$ws.Cells.Item(4, 13).Value = 29.801371
$ws.Cells.Item(4, 14).Value = 89.404113
$ws.Cells.Item(4, 15).Value = 0.08199584844219236
$ws.Cells.Item(4, 16).Value = 0.08199584844219235
$ws.Cells.Item(4, 17).Value = 355.5225282032173
$ws.Cells.Item(4, 18).Value = 3199.702753828956
$ws.Cells.Item(4, 19).Value = 0.01033338666764473
$ws.Cells.Item(4, 20).Value = 0.01033338666764473
$ws.Cells.Item(5, 7).Value = 11.92973733333333
$ws.Cells.Item(5, 8).Value = 35.789212
$ws.Cells.Item(5, 9).Value = 0.1260232909831994
$ws.Cells.Item(5, 10).Value = 0.1260232909831994
$ws.Cells.Item(5, 13).Value = 86.47679266666667
$ws.Cells.Item(5, 14).Value = 259.430378
$ws.Cells.Item(5, 15).Value = 0.2379332811655844
$ws.Cells.Item(5, 16).Value = 0.2379332811655844
$ws.Cells.Item(5, 17).Value = 1031.64542194246
$ws.Cells.Item(5, 18).Value = 9284.808797482136
$ws.Cells.Item(5, 19).Value = 0.02998513512691785
$ws.Cells.Item(5, 20).Value = 0.02998513512691785
$ws.Cells.Item(6, 9).Value = 0.7020418476259299
$ws.Cells.Item(6, 10).Value = 0.7020418476259298
$ws.Cells.Item(6, 13).Value = 3.795192333333334
$ws.Cells.Item(6, 14).Value = 11.385577
$ws.Cells.Item(6, 15).Value = 0.01044213755712683
$ws.Cells.Item(6, 16).Value = 0.01044213755712683
$ws.Cells.Item(6, 17).Value = 252.2184517799964
$ws.Cells.Item(6, 18).Value = 2269.966066019967
$ws.Cells.Item(6, 19).Value = 0.007330817543769436
$ws.Cells.Item(6, 20).Value = 0.007330817543769435
$ws.Cells.Item(7, 9).Value = 0.7020418476259299
$ws.Cells.Item(7, 10).Value = 0.7020418476259298
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.6696287328350964
$ws.Cells.Item(7, 16).Value = 0.6696287328350964
$ws.Cells.Item(7, 19).Value = 0.4701073928229613
$ws.Cells.Item(7, 20).Value = 0.4701073928229612
$ws.Cells.Item(8, 9).Value = 0.7020418476259299
$ws.Cells.Item(8, 10).Value = 0.7020418476259298
$ws.Cells.Item(8, 13).Value = 29.801371
$ws.Cells.Item(8, 14).Value = 89.404113
$ws.Cells.Item(8, 15).Value = 0.08199584844219236
$ws.Cells.Item(8, 16).Value = 0.08199584844219235
$ws.Cells.Item(8, 17).Value = 1980.520351636447
$ws.Cells.Item(8, 18).Value = 17824.68316472802
$ws.Cells.Item(8, 19).Value = 0.05756451693801246
$ws.Cells.Item(8, 20).Value = 0.05756451693801243
$ws.Cells.Item(9, 9).Value = 0.7020418476259299
$ws.Cells.Item(9, 10).Value = 0.7020418476259298
$ws.Cells.Item(9, 13).Value = 86.47679266666667
$ws.Cells.Item(9, 14).Value = 259.430378
$ws.Cells.Item(9, 15).Value = 0.2379332811655844
$ws.Cells.Item(9, 16).Value = 0.2379332811655844
$ws.Cells.Item(9, 17).Value = 5747.019082463649
$ws.Cells.Item(9, 18).Value = 51723.17174217285
$ws.Cells.Item(9, 19).Value = 0.1670391203211867
$ws.Cells.Item(9, 20).Value = 0.1670391203211867
$ws.Cells.Item(10, 7).Value = 16.00434166666667
$ws.Cells.Item(10, 8).Value = 48.013025
$ws.Cells.Item(10, 9).Value = 0.1690665729259037
$ws.Cells.Item(10, 10).Value = 0.1690665729259037
$ws.Cells.Item(10, 13).Value = 3.795192333333334
$ws.Cells.Item(10, 14).Value = 11.385577
$ws.Cells.Item(10, 15).Value = 0.01044213755712683
$ws.Cells.Item(10, 16).Value = 0.01044213755712683
$ws.Cells.Item(10, 17).Value = 60.73955479338056
$ws.Cells.Item(10, 18).Value = 546.655993140425
$ws.Cells.Item(10, 19).Value = 0.001765416410804302
$ws.Cells.Item(10, 20).Value = 0.001765416410804302
$ws.Cells.Item(11, 7).Value = 16.00434166666667
$ws.Cells.Item(11, 8).Value = 48.013025
$ws.Cells.Item(11, 9).Value = 0.1690665729259037
$ws.Cells.Item(11, 10).Value = 0.1690665729259037
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.6696287328350964
$ws.Cells.Item(11, 16).Value = 0.6696287328350964
$ws.Cells.Item(11, 17).Value = 3895.079037864211
$ws.Cells.Item(11, 18).Value = 35055.7113407779
$ws.Cells.Item(11, 19).Value = 0.1132118349931453
$ws.Cells.Item(11, 20).Value = 0.1132118349931453
$ws.Cells.Item(12, 7).Value = 16.00434166666667
$ws.Cells.Item(12, 8).Value = 48.013025
$ws.Cells.Item(12, 9).Value = 0.1690665729259037
$ws.Cells.Item(12, 10).Value = 0.1690665729259037
$ws.Cells.Item(12, 13).Value = 29.801371
$ws.Cells.Item(12, 14).Value = 89.404113
$ws.Cells.Item(12, 15).Value = 0.08199584844219236
$ws.Cells.Item(12, 16).Value = 0.08199584844219235
$ws.Cells.Item(12, 17).Value = 476.9513236190916
$ws.Cells.Item(12, 18).Value = 4292.561912571825
$ws.Cells.Item(12, 19).Value = 0.01386275709027327
$ws.Cells.Item(12, 20).Value = 0.01386275709027326
$ws.Cells.Item(13, 7).Value = 16.00434166666667
$ws.Cells.Item(13, 8).Value = 48.013025
$ws.Cells.Item(13, 9).Value = 0.1690665729259037
$ws.Cells.Item(13, 10).Value = 0.1690665729259037
$ws.Cells.Item(13, 13).Value = 86.47679266666667
$ws.Cells.Item(13, 14).Value = 259.430378
$ws.Cells.Item(13, 15).Value = 0.2379332811655844
$ws.Cells.Item(13, 16).Value = 0.2379332811655844
$ws.Cells.Item(13, 17).Value = 1384.004136074828
$ws.Cells.Item(13, 18).Value = 12456.03722467345
$ws.Cells.Item(13, 19).Value = 0.04022656443168083
$ws.Cells.Item(13, 20).Value = 0.04022656443168083
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.2715206666666667
$ws.Cells.Item(14, 8).Value = 0.814562
$ws.Cells.Item(14, 9).Value = 0.002868288464966955
$ws.Cells.Item(14, 10).Value = 0.002868288464966954
$ws.Cells.Item(14, 13).Value = 3.795192333333334
$ws.Cells.Item(14, 14).Value = 11.385577
$ws.Cells.Item(14, 15).Value = 0.01044213755712683
$ws.Cells.Item(14, 16).Value = 0.01044213755712683
$ws.Cells.Item(14, 17).Value = 1.030473152474889
$ws.Cells.Item(14, 18).Value = 9.274258372274002
$ws.Cells.Item(14, 19).Value = 0.00002995106270470511
$ws.Cells.Item(14, 20).Value = 0.00002995106270470511
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.2715206666666667
$ws.Cells.Item(15, 8).Value = 0.814562
$ws.Cells.Item(15, 9).Value = 0.002868288464966955
$ws.Cells.Item(15, 10).Value = 0.002868288464966954
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.6696287328350964
$ws.Cells.Item(15, 16).Value = 0.6696287328350964
$ws.Cells.Item(15, 17).Value = 66.08172201690579
$ws.Cells.Item(15, 18).Value = 594.7354981521521
$ws.Cells.Item(15, 19).Value = 0.001920688370201346
$ws.Cells.Item(15, 20).Value = 0.001920688370201345
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.2715206666666667
$ws.Cells.Item(16, 8).Value = 0.814562
$ws.Cells.Item(16, 9).Value = 0.002868288464966955
$ws.Cells.Item(16, 10).Value = 0.002868288464966954
$ws.Cells.Item(16, 13).Value = 29.801371
$ws.Cells.Item(16, 14).Value = 89.404113
$ws.Cells.Item(16, 15).Value = 0.08199584844219236
$ws.Cells.Item(16, 16).Value = 0.08199584844219235
$ws.Cells.Item(16, 17).Value = 8.091688121500667
$ws.Cells.Item(16, 18).Value = 72.825193093506
$ws.Cells.Item(16, 19).Value = 0.000235187746261919
$ws.Cells.Item(16, 20).Value = 0.0002351877462619189
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.2715206666666667
$ws.Cells.Item(17, 8).Value = 0.814562
$ws.Cells.Item(17, 9).Value = 0.002868288464966955
$ws.Cells.Item(17, 10).Value = 0.002868288464966954
$ws.Cells.Item(17, 13).Value = 86.47679266666667
$ws.Cells.Item(17, 14).Value = 259.430378
$ws.Cells.Item(17, 15).Value = 0.2379332811655844
$ws.Cells.Item(17, 16).Value = 0.2379332811655844
$ws.Cells.Item(17, 17).Value = 23.48023639604845
$ws.Cells.Item(17, 18).Value = 211.322127564436
$ws.Cells.Item(17, 19).Value = 0.0006824612857989848
$ws.Cells.Item(17, 20).Value = 0.0006824612857989847
